# registrar.xlsx - "reconfigurando e adicionando mais funcao"
#
# - Clear cell A2 (previously held the shared string "123testeregistro"),
#   which also drops that now-unused entry from the shared strings table.
# - Move/restore the active selection to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fully clear A2 (contents + formatting) so the cell element - and the
# shared string it referenced - disappear entirely, rather than leaving
# an empty placeholder cell behind.
$ws.Range("A2").Clear() | Out-Null

# Update the active selection to A2.
$ws.Range("A2").Select() | Out-Null
